# Commit: Explicitly directed to save to "Main Data" sheet.
#
# - Rename Sheet1 -> "Main Data"
# - Add a bold header row (A1:H1) with the tutor-center login log columns
# - Resize the columns to roughly match the new header content
# - Leave the selection spanning the first few data rows (A2:H4)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "Main Data"

# Header row content, left to right
$headers = @(
    "A-number",
    "Class Rank",
    "Major",
    "Course Prefix",
    "Course Name",
    "Date",
    "Day",
    "Time In"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Cells.Item(1, $i + 1)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
}

# Column widths tuned to the new headers
$ws.Columns.Item(1).ColumnWidth = 11.333333333333334
$ws.Columns.Item(2).ColumnWidth = 11.5
$ws.Columns.Item(3).ColumnWidth = 7.666666666666667
$ws.Columns.Item(4).ColumnWidth = 14.166666666666666
$ws.Columns.Item(5).ColumnWidth = 40.166666666666664
$ws.Columns.Item(6).ColumnWidth = 18.666666666666668
$ws.Columns.Item(7).ColumnWidth = 7.333333333333333
$ws.Columns.Item(8).ColumnWidth = 8.833333333333334

# Leave the same selection the author left the sheet in
[void]$ws.Range("A2:H4").Select()
